$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 5, 6, 8, 11, 12, 13)

$data = @()
$data += ,@(1.092944669801142, 0.04028668481693387, 0.346978685626695, 0.0842366385583162, 4.953339856236369, 0.07973214163530429, 0.617892686012965, 0.199028834499515, 0.2564379338705827)
$data += ,@(1.07909357416716, 0.03812766800085399, 0.3343197433336087, 0.08360673564669518, 4.74137407832751, 0.07973214163530429, 0.6045266538106091, 0.1968507443426901, 0.2534131723322197)
$data += ,@(1.071552234930095, 0.03677288371068954, 0.3264566223547121, 0.08323112675764044, 4.611407365189422, 0.07973214163530429, 0.5970327252605188, 0.1956125514579412, 0.2517483755631176)
$data += ,@(1.068720938829927, 0.03621330098867759, 0.3232291781494041, 0.0830808628611166, 4.558486500939665, 0.07973214163530429, 0.5941575649633393, 0.1951329034368641, 0.2511182923239197)
$data += ,@(1.068265401658834, 0.03611992494094451, 0.3226918529036027, 0.08305608061403014, 4.549701432159594, 0.07973214163530429, 0.5936909216755453, 0.1950547634398774, 0.2510165853798014)
$data += ,@(1.071513072278435, 0.03676536756071869, 0.3264131902129321, 0.0832290889158358, 4.610693492545323, 0.07973214163530429, 0.5969932271546412, 0.195605981845695, 0.2517396823995881)
$data += ,@(1.087968648099292, 0.039548200970561, 0.3426324133929768, 0.08401712970675668, 4.880214031325693, 0.07973214163530429, 0.6131357622472109, 0.1982572402479335, 0.2553550265628495)
$data += ,@(1.127902360285333, 0.04478189957535506, 0.3737398177964906, 0.0856513157949923, 5.410383309104475, 0.07973214163530429, 0.6504773667711845, 0.2042441855647468, 0.2639746959276152)
$data += ,@(1.161949957524769, 0.04850182590720209, 0.3961972016211632, 0.08690674287348443, 5.801217793267568, 0.07973214163530429, 0.6814253061023692, 0.2091253294510125, 0.2712461558002346)
$data += ,@(1.178469861633147, 0.05016929332177966, 0.4063336508878876, 0.08748991550165997, 5.979376596327995, 0.07973214163530429, 0.6962777568209049, 0.2114512260864387, 0.2747593367264187)
$data += ,@(1.184874437751716, 0.05079735209746161, 0.4101610952685348, 0.08771249296685824, 6.046898424745052, 0.07973214163530429, 0.7020140954311671, 0.2123471769292138, 0.2761193109635443)
$data += ,@(1.18348846919497, 0.05066223556797667, 0.4093372694447908, 0.0876644793148138, 6.032353788764851, 0.07973214163530429, 0.7007736776196793, 0.2121535421916718, 0.2758250980274823)
$data += ,@(1.178993783802838, 0.05022103067864236, 0.4066487552217097, 0.08750819211963901, 5.984930491512387, 0.07973214163530429, 0.6967474392268969, 0.2115246320541218, 0.2748706287570073)
$data += ,@(1.176260059253792, 0.04995034589659753, 0.4050005416525835, 0.08741268884477549, 5.955889898372448, 0.07973214163530429, 0.6942958644281418, 0.2111413845586867, 0.274289847109884)
$data += ,@(1.160891143136638, 0.04839237069224822, 0.3955331969748954, 0.08686887466291893, 5.789582431684863, 0.07973214163530429, 0.6804703016709936, 0.208975450213714, 0.271020700135189)
$data += ,@(1.151727363330593, 0.04743039546924166, 0.38970522187293, 0.08653835899194107, 5.687655025185109, 0.07973214163530429, 0.6721875383944678, 0.2076737405947426, 0.2690678362850747)
$data += ,@(1.146553655251751, 0.04687475252103468, 0.386345613582904, 0.08634939103523998, 5.629063288590118, 0.07973214163530429, 0.6674963274998618, 0.2069349547282116, 0.2679639249528734)
$data += ,@(1.144818580395594, 0.0466862138624009, 0.3852068049268809, 0.0862856046286602, 5.60923087938221, 0.07973214163530429, 0.6659204483468386, 0.2066865181753172, 0.2675934761638672)
$data += ,@(1.152692814740647, 0.04753304025391003, 0.3903263940486283, 0.08657342533816248, 5.69850179992423, 0.07973214163530429, 0.6730617113647384, 0.2078112826492173, 0.2692737214271759)
$data += ,@(1.18030993795017, 0.05035071337462682, 0.407438732385458, 0.08755405010843731, 5.99885828320464, 0.07973214163530429, 0.6979269964948287, 0.2117089459259347, 0.2751501752985703)
$data += ,@(1.199227172910156, 0.05217264008386735, 0.4185586659218075, 0.08820510580799379, 6.195492625782833, 0.07973214163530429, 0.7148312888899113, 0.2143448110729906, 0.2791633996231653)
$data += ,@(1.18905110359529, 0.05120197458720099, 0.4126294595659488, 0.08785669317104094, 6.090513180697087, 0.07973214163530429, 0.7057491203072459, 0.2129298932229773, 0.2770056444772706)
$data += ,@(1.152256039307588, 0.04748664263065905, 0.3900455902669648, 0.08655756857387331, 5.693597949439152, 0.07973214163530429, 0.6726662777221861, 0.2077490700347369, 0.2691805821507742)
$data += ,@(1.116275059608569, 0.0433888393403663, 0.3653962687009482, 0.08519964973835137, 5.266750071845024, 0.07973214163530429, 0.6397619673388135, 0.2025399877624707, 0.2614784560246122)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}

Write-Output "done"